# Update scraped 漫展 (con) listing stats: "想去人数" (F) and "最低票价" (G)
# counts on the 展览 (sheet1) and 全部类型 (sheet4) worksheets, matching a
# fresh data pull (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# ---- 展览 sheet ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = 50
$ws1.Range("F3").Value = 499
$ws1.Range("F5").Value = 1170
$ws1.Range("F6").Value = 14349
$ws1.Range("F7").Value = 16591
$ws1.Range("F9").Value = 102
$ws1.Range("F10").Value = 14
$ws1.Range("F12").Value = 202
$ws1.Range("F21").Value = 1266
$ws1.Range("F22").Value = 136
$ws1.Range("F23").Value = 71
$ws1.Range("F25").Value = 18
$ws1.Range("F27").Value = 6746
$ws1.Range("F29").Value = 20
$ws1.Range("F30").Value = 1121
$ws1.Range("F33").Value = 5756
$ws1.Range("F36").Value = 192
$ws1.Range("F37").Value = 4834

# ---- 全部类型 sheet ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = 50
$ws4.Range("F3").Value = 499
$ws4.Range("F5").Value = 1170
$ws4.Range("F6").Value = 14349
$ws4.Range("F7").Value = 16591
$ws4.Range("F9").Value = 102
$ws4.Range("F10").Value = 14
$ws4.Range("F12").Value = 202
$ws4.Range("F21").Value = 1266
$ws4.Range("F22").Value = 136
$ws4.Range("F23").Value = 71
$ws4.Range("F26").Value = 18
$ws4.Range("F28").Value = 6746
$ws4.Range("F30").Value = 20
$ws4.Range("F31").Value = 1121
$ws4.Range("F36").Value = 5756
$ws4.Range("F39").Value = 192
$ws4.Range("F40").Value = 4834
